# Thesis Charts and Graphs - add material properties table to Literature Review sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Literature Review")

# --- Remove the hidden chart-tracking defined names (_xlchart.v1.*) ---
while ($wb.Names.Count -gt 0) {
    $wb.Names.Item(1).Delete()
}

# --- Populate the new material-properties table (rows 19-22) ---
# Values are written in this particular order so that the shared-string
# table indices line up with how the sheet was actually authored.
$ws.Range("A20").Value = "PCL"
$ws.Range("B19").Value = "Tensile Modulus (Et)"
$ws.Range("C19").Value = "Compressive Modulus (Ecomp)"
$ws.Range("E19").Value = "Tensile Yield Strength"
$ws.Range("B20").Value = "440 ± 3 Mpa (ID46)`n417 ± 25 Mpa (ID133)"
$ws.Range("C20").Value = "455 ± 2 Mpa (ID46)"
$ws.Range("D20").Value = "414 ± 10 Mpa (ID46)"
$ws.Range("E20").Value = "17.82 ± 0.47 Mpa (ID46)`n14.7 ± 1.3 MPa (ID133)"
$ws.Range("A21").Value = "PLA"
$ws.Range("A22").Value = "PLCL (70/30)"
$ws.Range("B22").Value = "12 ± 1.2 Mpa (Secant Modulus at 0.2% strain) (ID31)"
$ws.Range("E22").Value = "17.2 ± 0.7 MPa (ID31)`n16.1 ± 3.2 Mpa (ID19)"
$ws.Range("E21").Value = "55.9 ± 6.5 (ID133)`n50 Mpa (ID434)"
$ws.Range("D19").Value = "Flexural Modulus (Eflex)"
$ws.Range("D21").Value = "3800 Mpa (ID343)"
$ws.Range("B21").Value = "3015 ± 86 Mpa (ID133)`n3600 Mpa (ID343)"
$ws.Range("A19").Value = "Material"

# --- Match formatting of the new table to the existing "inner cell" style ---
# (thin border on all sides, white fill) used elsewhere in the sheet, then
# center the text vertically, and wrap the multi-line measurement cells.
$ws.Range("B14").Copy()
$ws.Range("A19:E22").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A19:E22").VerticalAlignment = -4108  # xlCenter
$ws.Range("B20").WrapText = $true
$ws.Range("E20").WrapText = $true
$ws.Range("B21").WrapText = $true
$ws.Range("E21").WrapText = $true
$ws.Range("E22").WrapText = $true

# --- Row heights for the wrapped, multi-line rows ---
$ws.Rows.Item(20).RowHeight = 68
$ws.Rows.Item(21).RowHeight = 34
$ws.Rows.Item(22).RowHeight = 34

# --- Resize columns to fit the new, wider content ---
$ws.Columns.Item(1).ColumnWidth = 10.830729166666666
$ws.Columns.Item(2).ColumnWidth = 42.498697916666664
$ws.Columns.Item(3).ColumnWidth = 25.498697916666668
$ws.Columns.Item(4).ColumnWidth = 19.666666666666668
$ws.Columns.Item(5).ColumnWidth = 17.666666666666668

# --- Update selection / scroll position to reflect where the user ended up ---
$ws.Activate()
$ws.Range("B25").Select()

# --- Chart: normalize axis "crosses" to explicit autoZero on both axes ---
$cht = $ws.ChartObjects().Item(1).Chart
$cht.Axes(1).Crosses = "autoZero"
$cht.Axes(2).Crosses = "autoZero"
